$wb = $excel.ActiveWorkbook

$wsLog = $wb.Worksheets.Item("🗒更新日志")
$wsLog.Name = "🌳更新日志"

$wsYy = $wb.Worksheets.Item("yy的题目")
$wsYy.Name = "❓思考题"

$wsLog.Range("A8").NumberFormat = "@"
$wsLog.Range("A8").Value = "2023-04-13"
$wsLog.Range("A8").ClearFormats()
$wsLog.Range("B8").Value = "n by n 矩阵填入 -1 0 1，使得每行和为 0，每列和为 0，共有多少方法？"
